# Macroferia Regional de Talca - Kiwi
# Insert two new weekly price rows (dated 2022-08-22, serial 44795) right
# before the existing 2021-11-08 batch (row 231), shifting all subsequent
# rows down by 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows at position 231; this shifts rows 231:326 down to 233:328
# and brings along formatting (e.g. the date style on column D).
$ws.Rows("231:232").Insert()

# New row 231: Kiwi Hayward, Primera
$ws.Cells.Item(231, 1).Value = 5
$ws.Cells.Item(231, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(231, 3).Value = "Maule"
$ws.Cells.Item(231, 4).Value = 44795
$ws.Cells.Item(231, 5).Value = 7
$ws.Cells.Item(231, 6).Value = "Fruta"
$ws.Cells.Item(231, 7).Value = 100101
$ws.Cells.Item(231, 8).Value = "Berries"
$ws.Cells.Item(231, 9).Value = 100101007
$ws.Cells.Item(231, 10).Value = "Kiwi"
$ws.Cells.Item(231, 11).Value = "Hayward"
$ws.Cells.Item(231, 12).Value = "Primera"
$ws.Cells.Item(231, 13).Value = 250
$ws.Cells.Item(231, 14).Value = 7000
$ws.Cells.Item(231, 15).Value = 7000
$ws.Cells.Item(231, 16).Value = 7000
$ws.Cells.Item(231, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(231, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(231, 19).Value = 389
$ws.Cells.Item(231, 20).Value = 18

# New row 232: Kiwi Hayward, Segunda
$ws.Cells.Item(232, 1).Value = 5
$ws.Cells.Item(232, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(232, 3).Value = "Maule"
$ws.Cells.Item(232, 4).Value = 44795
$ws.Cells.Item(232, 5).Value = 7
$ws.Cells.Item(232, 6).Value = "Fruta"
$ws.Cells.Item(232, 7).Value = 100101
$ws.Cells.Item(232, 8).Value = "Berries"
$ws.Cells.Item(232, 9).Value = 100101007
$ws.Cells.Item(232, 10).Value = "Kiwi"
$ws.Cells.Item(232, 11).Value = "Hayward"
$ws.Cells.Item(232, 12).Value = "Segunda"
$ws.Cells.Item(232, 13).Value = 200
$ws.Cells.Item(232, 14).Value = 5000
$ws.Cells.Item(232, 15).Value = 5000
$ws.Cells.Item(232, 16).Value = 5000
$ws.Cells.Item(232, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(232, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(232, 19).Value = 278
$ws.Cells.Item(232, 20).Value = 18
